$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 586
$ws.Range("I2").Value = 1492
$ws.Range("J2").Value = 5904
$ws.Range("K2").Value = 29
$ws.Range("L2").Value = 1655
$ws.Range("M2").Value = 99
$ws.Range("N2").Value = 1057
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 19
$ws.Range("Q2").Value = 10
$ws.Range("R2").Value = 69
$ws.Range("S2").Value = 652
$ws.Range("T2").Value = 1105
$ws.Range("U2").Value = 64
$ws.Range("V2").Value = 9314
$ws.Range("W2").Value = 5
$ws.Range("X2").Value = 9147
$ws.Range("Y2").Value = 11
$ws.Range("Z2").Value = 122
$ws.Range("AA2").Value = 62
